$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-22 Sunday" "2024-09-23 Monday"

Replace-Text "695×6=4170" "212×2=424"
Replace-Text "281×8=2248" "470×4=1880"
Replace-Text "226×6=1356" "695×9=6255"
Replace-Text "520×9=4680" "358×8=2864"
Replace-Text "711×4=2844" "134×4=536"
Replace-Text "733×6=4398" "318×3=954"
Replace-Text "853×4=3412" "576×2=1152"
Replace-Text "734×8=5872" "867×8=6936"
Replace-Text "576×3=1728" "376×6=2256"
Replace-Text "426×8=3408" "858×6=5148"
Replace-Text "721×8=5768" "906×5=4530"
Replace-Text "990×9=8910" "709×4=2836"
Replace-Text "444×8=3552" "752×3=2256"
Replace-Text "371×8=2968" "974×3=2922"
Replace-Text "869×2=1738" "179×9=1611"
Replace-Text "884×9=7956" "157×6=942"
Replace-Text "842×5=4210" "272×8=2176"
Replace-Text "259×3=777" "735×9=6615"
Replace-Text "315×7=2205" "832×4=3328"
Replace-Text "750×7=5250" "964×6=5784"
Replace-Text "169×6=1014" "707×4=2828"
Replace-Text "656×4=2624" "879×6=5274"
Replace-Text "489×4=1956" "380×7=2660"
Replace-Text "928×5=4640" "712×8=5696"
Replace-Text "738×8=5904" "415×8=3320"
